# Apply the cryptos-list update described by the commit diff.
# Rows 2-34: price/volume refresh in place.
# Rows 35-51: a new "USDe" entry is inserted at row 35, shifting all
# subsequent coins down by one row (Bittensor falls off the bottom).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '65.329.32'
$ws.Range('E2').Value = '  +1.12%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.459.39'
$ws.Range('E3').Value = '  -0.42%  '

# Row 4
$ws.Range('E4').Value = '  +0.00%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '577.64'
$ws.Range('E5').Value = '  -0.65%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '161.37'
$ws.Range('E6').Value = '  +1.72%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.05%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.463.48'
$ws.Range('E8').Value = '  -0.36%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.581'
$ws.Range('E9').Value = '  +7.92%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.28'
$ws.Range('E10').Value = '  -4.43%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.126'
$ws.Range('E11').Value = '  +0.20%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.440'
$ws.Range('E12').Value = '  -0.44%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.057.04'
$ws.Range('E13').Value = '  -0.43%  '

# Row 14
$ws.Range('E14').Value = '  -2.12%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.0000195'
$ws.Range('E15').Value = '  -1.57%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '28.38'
$ws.Range('E16').Value = '  +2.26%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '65.270.48'
$ws.Range('E17').Value = '  +0.98%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.464.77'
$ws.Range('E18').Value = '  -0.64%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.37'
$ws.Range('E19').Value = '  -1.36%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.25'
$ws.Range('E20').Value = '  -0.78%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '388.25'
$ws.Range('E21').Value = '  -2.59%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '8.22'
$ws.Range('E22').Value = '  -4.05%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '73.42'
$ws.Range('E23').Value = '  +1.67%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.545'
$ws.Range('E24').Value = '  -0.10%  '

# Row 25
$ws.Range('E25').Value = '  +0.46%  '

# Row 26
$ws.Range('E26').Value = '  +3.55%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.78'
$ws.Range('E27').Value = '  +0.55%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.181'
$ws.Range('E28').Value = '  -0.47%  '

# Row 29
$ws.Range('E29').Value = '  -0.15%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.31'
$ws.Range('E30').Value = '  +6.16%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.45'
$ws.Range('E31').Value = '  +3.48%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.05'
$ws.Range('E32').Value = '  -0.45%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.63'
$ws.Range('E33').Value = '  -0.22%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '23.71'
$ws.Range('E34').Value = '  -0.35%  '

# Row 35
$ws.Range('B35').Value = 'USDe'
$ws.Range('C35').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').Value = '  +0.02%  '

# Row 36
$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '7.11'
$ws.Range('E36').Value = '  +2.44%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.50'
$ws.Range('E37').Value = '  -0.13%  '

# Row 38
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '162.84'
$ws.Range('E38').Value = '  +2.43%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.93'
$ws.Range('E39').Value = '  +1.79%  '

# Row 40
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.040.80'
$ws.Range('E40').Value = '  +4.22%  '

# Row 41
$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0770'
$ws.Range('E41').Value = '  -0.84%  '

# Row 42
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '27.21'
$ws.Range('E42').Value = '  -3.58%  '

# Row 43
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '4.57'
$ws.Range('E43').Value = '  +3.00%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0319'
$ws.Range('E44').Value = '  -1.01%  '

# Row 45
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '42.64'
$ws.Range('E45').Value = '  +2.96%  '

# Row 46
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.772'
$ws.Range('E46').Value = '  -1.08%  '

# Row 47
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '25.05'
$ws.Range('E47').Value = '  +8.76%  '

# Row 48
$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.10'
$ws.Range('E48').Value = '  -2.23%  '

# Row 49
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.21'

# Row 50
$ws.Range('B50').Value = 'SuiNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.871'
$ws.Range('E50').Value = '  +3.58%  '

# Row 51
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '6.64'
$ws.Range('E51').Value = '  +2.26%  '
